# 141: 31/12 21:42 LP1912+6203+6173 — append newest scrape rows to the
# three schedule sheets (LP1912, LP1912-215, 6203-6173) and refresh the
# "Ultima actualizacion" / "Total filas" header cells on each sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912": columns A(meta) B=Hora_Scrap C=Hora_Llegada D=Linea
#                 E=Minutos F=Parada G=Fecha
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 31/12/2025 18:42:12"
$ws1.Range("A3").Value = "Total filas: 1238"

$sheet1Rows = @()
$sheet1Rows += ,@(1222, "18:42:01", "18:52", "15_ABASTO", 10, "LP1912", "31/12/2025")
$sheet1Rows += ,@(1223, "18:42:01", "18:54", "16_SANTA ANA", 12, "LP1912", "31/12/2025")
$sheet1Rows += ,@(1224, "18:42:01", "19:00", "10_OLMOS", 18, "LP1912", "31/12/2025")
$sheet1Rows += ,@(1225, "18:42:01", "19:02", "17_ROMERO", 20, "LP1912", "31/12/2025")
$sheet1Rows += ,@(1226, "18:42:01", "19:04", "23_HERNANDEZ", 22, "LP1912", "31/12/2025")
$sheet1Rows += ,@(1227, "18:42:01", "19:06", "16_SANTA ANA", 24, "LP1912", "31/12/2025")
$sheet1Rows += ,@(1228, "18:42:01", "19:15", "14_ABASTO", 33, "LP1912", "31/12/2025")
$sheet1Rows += ,@(1229, "18:42:01", "19:18", "16_SANTA ANA", 36, "LP1912", "31/12/2025")
$sheet1Rows += ,@(1230, "18:42:01", "19:22", "215C_EL PATO", 40, "LP1912", "31/12/2025")
$sheet1Rows += ,@(1231, "18:42:01", "19:32", "215_EL PELIGRO", 50, "LP1912", "31/12/2025")
$sheet1Rows += ,@(1232, "18:42:01", "19:34", "23_HERNANDEZ", 52, "LP1912", "31/12/2025")
$sheet1Rows += ,@(1233, "18:42:01", "19:41", "17X38_ROMERO", 59, "LP1912", "31/12/2025")
$sheet1Rows += ,@(1234, "18:42:01", "19:45", "11_ETCHEVERRY", 63, "LP1912", "31/12/2025")
$sheet1Rows += ,@(1235, "18:42:01", "19:52", "81_EL PELIGRO", 70, "LP1912", "31/12/2025")
$sheet1Rows += ,@(1236, "18:42:01", "19:59", "14X44_ABASTO", 77, "LP1912", "31/12/2025")
$sheet1Rows += ,@(1237, "18:42:01", "20:02", "215C_EL PATO", 80, "LP1912", "31/12/2025")
$sheet1Rows += ,@(1238, "18:42:01", "20:11", "23_HERNANDEZ", 89, "LP1912", "31/12/2025")
$sheet1Rows += ,@(1239, "18:42:01", "20:15", "11_ETCHEVERRY", 93, "LP1912", "31/12/2025")

foreach ($row in $sheet1Rows) {
    $r = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $ws1.Cells.Item($r, 7).Value = $row[6]
}

# ---------------------------------------------------------------------------
# Sheet "LP1912-215": columns A(meta) B=Fecha C=Hora_Scrap D=Hora_Llegada
#                      E=Linea F=Minutos G=Parada
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 31/12/2025 18:42:12"
$ws2.Range("A3").Value = "Total filas: 86"

$sheet2Rows = @()
$sheet2Rows += ,@(85, "31/12/2025", "18:42:01", "19:22", "215C_EL PATO", 40, "LP1912")
$sheet2Rows += ,@(86, "31/12/2025", "18:42:01", "19:32", "215_EL PELIGRO", 50, "LP1912")
$sheet2Rows += ,@(87, "31/12/2025", "18:42:01", "20:02", "215C_EL PATO", 80, "LP1912")

foreach ($row in $sheet2Rows) {
    $r = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $ws2.Cells.Item($r, 6).Value = $row[5]
    $ws2.Cells.Item($r, 7).Value = $row[6]
}

# ---------------------------------------------------------------------------
# Sheet "6203-6173": columns A(meta) B=Fecha C=Hora_Scrap D=Hora_Llegada
#                     E=Linea F=Minutos G=Parada
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 31/12/2025 18:42:12"
$ws3.Range("A3").Value = "Total filas: 145"

$sheet3Rows = @()
$sheet3Rows += ,@(146, "31/12/2025", "18:42:12", "19:11", "215B_LP-P MOR-1 Y 57", 29, "L6173")

foreach ($row in $sheet3Rows) {
    $r = $row[0]
    $ws3.Cells.Item($r, 2).Value = $row[1]
    $ws3.Cells.Item($r, 3).Value = $row[2]
    $ws3.Cells.Item($r, 4).Value = $row[3]
    $ws3.Cells.Item($r, 5).Value = $row[4]
    $ws3.Cells.Item($r, 6).Value = $row[5]
    $ws3.Cells.Item($r, 7).Value = $row[6]
}
